$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.2384341637010676
$ws.Cells.Item(2, 3).Value = 0.4697508896797153
$ws.Cells.Item(2, 10).Value = 0.01779359430604982
$ws.Cells.Item(2, 16).Value = 0.1779359430604982
$ws.Cells.Item(2, 19).Value = 0.09608540925266904
$ws.Cells.Item(3, 10).Value = 0.04477611940298507
$ws.Cells.Item(3, 16).Value = 0.7238805970149254
$ws.Cells.Item(3, 19).Value = 0.2313432835820896
$ws.Cells.Item(4, 10).Value = 0.03571428571428571
$ws.Cells.Item(4, 16).Value = 0.6785714285714286
$ws.Cells.Item(4, 19).Value = 0.2857142857142857
$ws.Cells.Item(6, 2).Value = 0.08016877637130802
$ws.Cells.Item(6, 4).Value = 0.008438818565400843
$ws.Cells.Item(6, 6).Value = 0.09282700421940929
$ws.Cells.Item(6, 10).Value = 0.1940928270042194
$ws.Cells.Item(6, 15).Value = 0.01265822784810127
$ws.Cells.Item(6, 17).Value = 0.2151898734177215
$ws.Cells.Item(6, 18).Value = 0.04641350210970464
$ws.Cells.Item(6, 19).Value = 0.350210970464135
$ws.Cells.Item(7, 2).Value = 0.09090909090909091
$ws.Cells.Item(7, 4).Value = 0.0374331550802139
$ws.Cells.Item(7, 6).Value = 0.106951871657754
$ws.Cells.Item(7, 10).Value = 0.1229946524064171
$ws.Cells.Item(7, 17).Value = 0.1711229946524064
$ws.Cells.Item(7, 18).Value = 0.1283422459893048
$ws.Cells.Item(7, 19).Value = 0.3422459893048128
$ws.Cells.Item(8, 2).Value = 0.07175925925925926
$ws.Cells.Item(8, 4).Value = 0.009259259259259259
$ws.Cells.Item(8, 5).Value = 0.002314814814814815
$ws.Cells.Item(8, 6).Value = 0.05092592592592592
$ws.Cells.Item(8, 10).Value = 0.1226851851851852
$ws.Cells.Item(8, 15).Value = 0.01157407407407407
$ws.Cells.Item(8, 17).Value = 0.150462962962963
$ws.Cells.Item(8, 18).Value = 0.1388888888888889
$ws.Cells.Item(8, 19).Value = 0.4421296296296297
$ws.Cells.Item(9, 2).Value = 0.08389261744966443
$ws.Cells.Item(9, 4).Value = 0.02013422818791946
$ws.Cells.Item(9, 5).Value = 0.003355704697986577
$ws.Cells.Item(9, 6).Value = 0.06375838926174497
$ws.Cells.Item(9, 10).Value = 0.1241610738255034
$ws.Cells.Item(9, 15).Value = 0.02013422818791946
$ws.Cells.Item(9, 17).Value = 0.1644295302013423
$ws.Cells.Item(9, 18).Value = 0.1006711409395973
$ws.Cells.Item(9, 19).Value = 0.4194630872483222
$ws.Cells.Item(10, 2).Value = 0.09501557632398754
$ws.Cells.Item(10, 4).Value = 0.02725856697819315
$ws.Cells.Item(10, 6).Value = 0.06853582554517133
$ws.Cells.Item(10, 10).Value = 0.1168224299065421
$ws.Cells.Item(10, 15).Value = 0.0132398753894081
$ws.Cells.Item(10, 17).Value = 0.1853582554517134
$ws.Cells.Item(10, 18).Value = 0.1004672897196262
$ws.Cells.Item(10, 19).Value = 0.3933021806853583
$ws.Cells.Item(11, 7).Value = 0.1362007168458781
$ws.Cells.Item(11, 10).Value = 0.1003584229390681
$ws.Cells.Item(11, 11).Value = 0.1863799283154122
$ws.Cells.Item(11, 12).Value = 0.5663082437275986
$ws.Cells.Item(11, 19).Value = 0.01075268817204301
$ws.Cells.Item(12, 7).Value = 0.8343558282208589
$ws.Cells.Item(12, 10).Value = 0.09202453987730061
$ws.Cells.Item(12, 11).Value = 0.006134969325153374
$ws.Cells.Item(12, 12).Value = 0.03067484662576687
$ws.Cells.Item(12, 19).Value = 0.03680981595092025
$ws.Cells.Item(13, 7).Value = 0.6285714285714286
$ws.Cells.Item(13, 10).Value = 0.2571428571428571
$ws.Cells.Item(13, 19).Value = 0.1142857142857143
$ws.Cells.Item(15, 6).Value = 0.01785714285714286
$ws.Cells.Item(15, 8).Value = 0.1785714285714286
$ws.Cells.Item(15, 9).Value = 0.0625
$ws.Cells.Item(15, 10).Value = 0.3839285714285715
$ws.Cells.Item(15, 11).Value = 0.05803571428571429
$ws.Cells.Item(15, 13).Value = 0.01339285714285714
$ws.Cells.Item(15, 15).Value = 0.1071428571428571
$ws.Cells.Item(15, 19).Value = 0.1785714285714286
$ws.Cells.Item(16, 6).Value = 0.005649717514124294
$ws.Cells.Item(16, 8).Value = 0.1525423728813559
$ws.Cells.Item(16, 9).Value = 0.1638418079096045
$ws.Cells.Item(16, 10).Value = 0.3898305084745763
$ws.Cells.Item(16, 11).Value = 0.1073446327683616
$ws.Cells.Item(16, 13).Value = 0.01129943502824859
$ws.Cells.Item(16, 14).Value = 0.005649717514124294
$ws.Cells.Item(16, 15).Value = 0.04519774011299435
$ws.Cells.Item(16, 19).Value = 0.1186440677966102
$ws.Cells.Item(17, 6).Value = 0.01624129930394431
$ws.Cells.Item(17, 8).Value = 0.1624129930394431
$ws.Cells.Item(17, 9).Value = 0.1554524361948956
$ws.Cells.Item(17, 10).Value = 0.4129930394431555
$ws.Cells.Item(17, 11).Value = 0.1020881670533643
$ws.Cells.Item(17, 13).Value = 0.009280742459396751
$ws.Cells.Item(17, 14).Value = 0.002320185614849188
$ws.Cells.Item(17, 15).Value = 0.07192575406032482
$ws.Cells.Item(17, 19).Value = 0.06728538283062645
$ws.Cells.Item(18, 6).Value = 0.02371541501976284
$ws.Cells.Item(18, 8).Value = 0.150197628458498
$ws.Cells.Item(18, 9).Value = 0.1146245059288538
$ws.Cells.Item(18, 10).Value = 0.4150197628458498
$ws.Cells.Item(18, 11).Value = 0.08300395256916997
$ws.Cells.Item(18, 13).Value = 0.01185770750988142
$ws.Cells.Item(18, 14).Value = 0.01185770750988142
$ws.Cells.Item(18, 15).Value = 0.06719367588932806
$ws.Cells.Item(18, 19).Value = 0.1225296442687747
$ws.Cells.Item(19, 6).Value = 0.01452599388379205
$ws.Cells.Item(19, 8).Value = 0.2010703363914373
$ws.Cells.Item(19, 9).Value = 0.1230886850152905
$ws.Cells.Item(19, 10).Value = 0.3654434250764526
$ws.Cells.Item(19, 11).Value = 0.09709480122324159
$ws.Cells.Item(19, 13).Value = 0.01758409785932722
$ws.Cells.Item(19, 14).Value = 0.001529051987767584
$ws.Cells.Item(19, 15).Value = 0.06651376146788991
$ws.Cells.Item(19, 19).Value = 0.1131498470948012
